# Stundenliste-AlexanderSchmid.xlsx
# "Potentially fixed error of channels not appearing before measuredData is received"
#
# Two more work-log entries were added to the "Tätigkeit" sheet for row 28/29
# (8h each), which in turn bumps the running "Zwischensumme" total carried
# down through the remaining rows. The active selection also ended up on
# D29 (the last edited cell), and the date column's number format was
# switched from the US m/d/yyyy style to the European D/M/YYYY style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: new day worked, 8 hours, with a description of the work done.
$ws.Cells.Item(28, 1).Value = 44443
$ws.Cells.Item(28, 2).Value = 8
$ws.Cells.Item(28, 4).Value = "Weitere Implementierung der Anzeige der MeasuredData"

# Row 29: another new day worked, 8 hours, with its own description.
$ws.Cells.Item(29, 1).Value = 44444
$ws.Cells.Item(29, 2).Value = 8
$ws.Cells.Item(29, 4).Value = "Bugfixing der Anzeige der MeasuredData"

# Switch the "Datum" column from US-style to European-style date formatting.
$ws.Range("A5:A26").NumberFormat = "D/M/YYYY"
$ws.Range("A27:A38").NumberFormat = "DD/MM/YY"

# Leave the selection where the author's last edit was (D29).
$ws.Range("D29").Select()
